# Scheduled-runner refresh: update cached market-board stats (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H:N) for the
# leve rows whose source prices changed, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1626.2106
$ws.Range("I129").Value = 726.5333000000001
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 2179.5999
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 2820.4001
$ws.Range("N129").Value = -25000
$ws.Range("H132").Value = 3266.776
$ws.Range("I132").Value = 3469.52
$ws.Range("K132").Value = 10408.56
$ws.Range("M132").Value = -7878.559999999999
$ws.Range("H137").Value = 895660.75
$ws.Range("I137").Value = 2274317.5
$ws.Range("K137").Value = 6822952.5
$ws.Range("M137").Value = -6820402.5
$ws.Range("H138").Value = 3203.868
$ws.Range("J138").Value = 2904.5757
$ws.Range("L138").Value = 8713.7271
$ws.Range("N138").Value = -18993.7271
$ws.Range("H141").Value = 3076.8235
$ws.Range("I141").Value = 3063.3635
$ws.Range("J141").Value = 3101.5
$ws.Range("K141").Value = 9190.0905
$ws.Range("L141").Value = 9304.5
$ws.Range("M141").Value = -4010.0905
$ws.Range("N141").Value = -19664.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1567
$ws.Range("I2").Value = 678.8570999999999
$ws.Range("K2").Value = 678.8570999999999
$ws.Range("M2").Value = -565.8570999999999
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H116").Value = 1567
$ws.Range("I116").Value = 678.8570999999999
$ws.Range("K116").Value = 678.8570999999999
$ws.Range("M116").Value = 1615.1429
$ws.Range("H122").Value = 4484.7144
$ws.Range("I122").Value = 1863
$ws.Range("K122").Value = 5589
$ws.Range("M122").Value = -3139

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1567
$ws.Range("I3").Value = 678.8570999999999
$ws.Range("K3").Value = 678.8570999999999
$ws.Range("M3").Value = -564.8570999999999
$ws.Range("H64").Value = 881.9167
$ws.Range("I64").Value = 872.25
$ws.Range("K64").Value = 872.25
$ws.Range("M64").Value = -647.25
$ws.Range("H67").Value = 881.9167
$ws.Range("I67").Value = 872.25
$ws.Range("K67").Value = 872.25
$ws.Range("M67").Value = -92.25
$ws.Range("H86").Value = 2456.1333
$ws.Range("J86").Value = 2583.7144
$ws.Range("L86").Value = 2583.7144
$ws.Range("N86").Value = -4829.7144
$ws.Range("H89").Value = 2456.1333
$ws.Range("J89").Value = 2583.7144
$ws.Range("L89").Value = 12918.572
$ws.Range("N89").Value = -24150.572
$ws.Range("H94").Value = 71429710
$ws.Range("I94").Value = 95238440
$ws.Range("K94").Value = 95238440
$ws.Range("M94").Value = -95237989

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1809.1111
$ws.Range("I22").Value = 1785.25
$ws.Range("K22").Value = 1785.25
$ws.Range("M22").Value = -1435.25
$ws.Range("H31").Value = 4052.5
$ws.Range("I31").Value = 2715.9167
$ws.Range("J31").Value = 8062.25
$ws.Range("K31").Value = 2715.9167
$ws.Range("L31").Value = 8062.25
$ws.Range("M31").Value = -2420.9167
$ws.Range("N31").Value = -8652.25
$ws.Range("H34").Value = 4052.5
$ws.Range("I34").Value = 2715.9167
$ws.Range("J34").Value = 8062.25
$ws.Range("K34").Value = 2715.9167
$ws.Range("L34").Value = 8062.25
$ws.Range("M34").Value = -2513.9167
$ws.Range("N34").Value = -8466.25
$ws.Range("H58").Value = 1976.0741
$ws.Range("I58").Value = 998.6667
$ws.Range("K58").Value = 998.6667
$ws.Range("M58").Value = -795.6667
$ws.Range("H122").Value = 4540.778
$ws.Range("J122").Value = 5380.923
$ws.Range("L122").Value = 16142.769
$ws.Range("N122").Value = -21042.769
$ws.Range("H136").Value = 1976.0741
$ws.Range("I136").Value = 998.6667
$ws.Range("K136").Value = 2996.0001
$ws.Range("M136").Value = -446.0001000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5961.9395
$ws.Range("I131").Value = 15530.7
$ws.Range("J131").Value = 1801.6086
$ws.Range("K131").Value = 46592.10000000001
$ws.Range("L131").Value = 5404.825800000001
$ws.Range("M131").Value = -41552.10000000001
$ws.Range("N131").Value = -15484.8258
$ws.Range("H140").Value = 15223.186
$ws.Range("I140").Value = 11825.059
$ws.Range("K140").Value = 35475.177
$ws.Range("M140").Value = -30295.177

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 43479600
$ws.Range("J80").Value = 1891.5
$ws.Range("L80").Value = 1891.5
$ws.Range("N80").Value = -3887.5
$ws.Range("H83").Value = 43479600
$ws.Range("J83").Value = 1891.5
$ws.Range("L83").Value = 9457.5
$ws.Range("N83").Value = -19441.5
$ws.Range("H102").Value = 1227.1666
$ws.Range("I102").Value = 1068.0625
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1068.0625
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 553.9375
$ws.Range("N102").Value = -5744

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H9").Value = 2320
$ws.Range("J9").Value = 5500
$ws.Range("L9").Value = 5500
$ws.Range("N9").Value = -5948
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H100").Value = 2759.5
$ws.Range("I100").Value = 1370.8572
$ws.Range("J100").Value = 5999.6665
$ws.Range("K100").Value = 1370.8572
$ws.Range("L100").Value = 5999.6665
$ws.Range("M100").Value = -829.8571999999999
$ws.Range("N100").Value = -7081.6665
$ws.Range("H122").Value = 10372.417
$ws.Range("I122").Value = 9647.200000000001
$ws.Range("K122").Value = 28941.6
$ws.Range("M122").Value = -26491.6
$ws.Range("H136").Value = 3899.8333
$ws.Range("I136").Value = 4014.2144
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 12042.6432
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -9492.643199999999
$ws.Range("N136").Value = -15598.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1446.5
$ws.Range("I23").Value = 1446.5
$ws.Range("K23").Value = 1446.5
$ws.Range("M23").Value = -1217.5
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H62").Value = 12667.667
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 12667.667
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880
$ws.Range("H81").Value = 3797.5
$ws.Range("I81").Value = 3557
$ws.Range("K81").Value = 7114
$ws.Range("M81").Value = -6053
$ws.Range("H84").Value = 3797.5
$ws.Range("I84").Value = 3557
$ws.Range("K84").Value = 35570
$ws.Range("M84").Value = -30266
$ws.Range("H96").Value = 3195.4
$ws.Range("I96").Value = 1989.5
$ws.Range("K96").Value = 1989.5
$ws.Range("M96").Value = -616.5
$ws.Range("H136").Value = 1280.7273
$ws.Range("I136").Value = 1029.2222
$ws.Range("J136").Value = 2412.5
$ws.Range("K136").Value = 3087.6666
$ws.Range("L136").Value = 7237.5
$ws.Range("M136").Value = -537.6665999999996
$ws.Range("N136").Value = -12337.5
